# Shift Rota Generation - Bug fixing with Lists
# Updates the Cal_Primary (D), Cal_Standby (E) and BAS_FinC (F) columns
# of the rota table on Sheet1 to correct the rotation assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Divik"

$ws.Range("D4").Value = "Sushvin"
$ws.Range("F4").Value = "Naveen"

$ws.Range("E5").Value = "Kapil"
$ws.Range("F5").Value = "Sushvin"

$ws.Range("D6").Value = "Kapil"
$ws.Range("E6").Value = "Naveen"

$ws.Range("D7").Value = "Sushvin"
$ws.Range("E7").Value = "Divik"
$ws.Range("F7").Value = "Naveen"

$ws.Range("D8").Value = "Divik"
$ws.Range("E8").Value = "Naveen"
$ws.Range("F8").Value = "Kapil"

$ws.Range("D9").Value = "Kapil"
$ws.Range("E9").Value = "Naveen"
$ws.Range("F9").Value = "Sushvin"

$ws.Range("D10").Value = "Naveen"
$ws.Range("F10").Value = "Kapil"

$ws.Range("D12").Value = "Sushvin"
$ws.Range("E12").Value = "Kapil"

$ws.Range("D13").Value = "Kapil"
$ws.Range("E13").Value = "Divik"
$ws.Range("F13").Value = "Naveen"

$ws.Range("D14").Value = "Naveen"
$ws.Range("E14").Value = "Kapil"
$ws.Range("F14").Value = "Sushvin"

$ws.Range("F15").Value = "Sushvin"

$ws.Range("D16").Value = "Naveen"
$ws.Range("E16").Value = "Sushvin"
$ws.Range("F16").Value = "Divik"

$ws.Range("D17").Value = "Naveen"
$ws.Range("E17").Value = "Sushvin"
$ws.Range("F17").Value = "Divik"

$ws.Range("E18").Value = "Sushvin"
$ws.Range("F18").Value = "Naveen"

$ws.Range("D19").Value = "Sushvin"

$ws.Range("D20").Value = "Naveen"
$ws.Range("E20").Value = "Divik"
$ws.Range("F20").Value = "Kapil"

$ws.Range("D21").Value = "Kapil"
$ws.Range("E21").Value = "Sushvin"
$ws.Range("F21").Value = "Divik"

$ws.Range("E22").Value = "Divik"
$ws.Range("F22").Value = "Naveen"
